$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.720.22'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '1.846.97'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'313.52"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = "'0.4305"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.53%  '
$ws.Range('D8').Value = "'0.3654"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = "'45.04"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('D10').Value = "'0.07342"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('D11').Value = "'0.8772"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.59%  '
$ws.Range('D12').Value = "'20.77"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').Value = '1.769.01'
$ws.Range('E13').Value = '  -4.09%  '
$ws.Range('D14').Value = "'5.337"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').Value = "'6.521"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = "'0.06929"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.34%  '
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').Value = "'79.94"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.27%  '
$ws.Range('D19').Value = "'0.000008995"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('D20').Value = "'0.9995"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').Value = "'15.35"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.87%  '
$ws.Range('D22').Value = '27.566.16'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = "'4.975"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = "'10.38"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.12%  '
$ws.Range('D25').Value = '1.994.14'
$ws.Range('E25').Value = '  -2.52%  '
$ws.Range('E26').Value = '  -3.18%  '
$ws.Range('D27').Value = "'156.13"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.31%  '
$ws.Range('D28').Value = "'18.63"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.18%  '
$ws.Range('D29').Value = "'120.21"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.67%  '
$ws.Range('D30').Value = "'5.250"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.63%  '
$ws.Range('D31').Value = "'1.863"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.84%  '
$ws.Range('D32').Value = "'0.08904"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('D33').Value = "'0.7544"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.20%  '
$ws.Range('D34').Value = "'4.543"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('D35').Value = "'2.970"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('D36').Value = "'1.122"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.91%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = "'1.109"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.84%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = "'0.05411"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.01930"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = "'2.837"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = "'0.5083"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = "'0.1657"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.83%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = "'6.658"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.22%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = "'8.346"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.23%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = "'0.06535"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.47%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'10.31"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = "'0.4667"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.18%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'104.46"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = "'0.9998"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = "'1.622"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.09%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = "'64.03"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.33%  '
